$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 100.8744713333333
$ws.Cells.Item(2, 8).Value = 302.623414
$ws.Cells.Item(2, 9).Value = 0.1452075237922473
$ws.Cells.Item(2, 10).Value = 0.1452075237922473
$ws.Cells.Item(2, 13).Value = 9.423852333333334
$ws.Cells.Item(2, 14).Value = 28.271557
$ws.Cells.Item(2, 15).Value = 0.06654336290212845
$ws.Cells.Item(2, 16).Value = 0.06654336290212845
$ws.Cells.Item(2, 17).Value = 950.6261220484
$ws.Cells.Item(2, 18).Value = 8555.6350984356
$ws.Cells.Item(2, 19).Value = 0.009662596951826965
$ws.Cells.Item(2, 20).Value = 0.009662596951826965
$ws.Cells.Item(3, 7).Value = 100.8744713333333
$ws.Cells.Item(3, 8).Value = 302.623414
$ws.Cells.Item(3, 9).Value = 0.1452075237922473
$ws.Cells.Item(3, 10).Value = 0.1452075237922473
$ws.Cells.Item(3, 15).Value = 0.3572423751649123
$ws.Cells.Item(3, 16).Value = 0.3572423751649123
$ws.Cells.Item(3, 17).Value = 5103.49821414748
$ws.Cells.Item(3, 18).Value = 45931.48392732732
$ws.Cells.Item(3, 19).Value = 0.05187428069135795
$ws.Cells.Item(3, 20).Value = 0.05187428069135795
$ws.Cells.Item(4, 7).Value = 100.8744713333333
$ws.Cells.Item(4, 8).Value = 302.623414
$ws.Cells.Item(4, 9).Value = 0.1452075237922473
$ws.Cells.Item(4, 10).Value = 0.1452075237922473
$ws.Cells.Item(4, 13).Value = 26.84076266666667
$ws.Cells.Item(4, 14).Value = 80.522288
$ws.Cells.Item(4, 15).Value = 0.1895270158659356
$ws.Cells.Item(4, 16).Value = 0.1895270158659356
$ws.Cells.Item(4, 17).Value = 2707.547744183471
$ws.Cells.Item(4, 18).Value = 24367.92969765123
$ws.Cells.Item(4, 19).Value = 0.02752074866562648
$ws.Cells.Item(4, 20).Value = 0.02752074866562648
$ws.Cells.Item(5, 7).Value = 100.8744713333333
$ws.Cells.Item(5, 8).Value = 302.623414
$ws.Cells.Item(5, 9).Value = 0.1452075237922473
$ws.Cells.Item(5, 10).Value = 0.1452075237922473
$ws.Cells.Item(5, 13).Value = 54.762539
$ws.Cells.Item(5, 14).Value = 164.287617
$ws.Cells.Item(5, 15).Value = 0.3866872460670236
$ws.Cells.Item(5, 16).Value = 0.3866872460670236
$ws.Cells.Item(5, 17).Value = 5524.142170496049
$ws.Cells.Item(5, 18).Value = 49717.27953446444
$ws.Cells.Item(5, 19).Value = 0.05614989748343593
$ws.Cells.Item(5, 20).Value = 0.05614989748343593
$ws.Cells.Item(6, 9).Value = 0.7769829249672668
$ws.Cells.Item(6, 10).Value = 0.776982924967267
$ws.Cells.Item(6, 13).Value = 9.423852333333334
$ws.Cells.Item(6, 14).Value = 28.271557
$ws.Cells.Item(6, 15).Value = 0.06654336290212845
$ws.Cells.Item(6, 16).Value = 0.06654336290212845
$ws.Cells.Item(6, 17).Value = 5086.652850827631
$ws.Cells.Item(6, 18).Value = 45779.87565744867
$ws.Cells.Item(6, 19).Value = 0.05170305674485408
$ws.Cells.Item(6, 20).Value = 0.05170305674485409
$ws.Cells.Item(7, 9).Value = 0.7769829249672668
$ws.Cells.Item(7, 10).Value = 0.776982924967267
$ws.Cells.Item(7, 15).Value = 0.3572423751649123
$ws.Cells.Item(7, 16).Value = 0.3572423751649123
$ws.Cells.Item(7, 19).Value = 0.2775712255778873
$ws.Cells.Item(7, 20).Value = 0.2775712255778873
$ws.Cells.Item(8, 9).Value = 0.7769829249672668
$ws.Cells.Item(8, 10).Value = 0.776982924967267
$ws.Cells.Item(8, 13).Value = 26.84076266666667
$ws.Cells.Item(8, 14).Value = 80.522288
$ws.Cells.Item(8, 15).Value = 0.1895270158659356
$ws.Cells.Item(8, 16).Value = 0.1895270158659356
$ws.Cells.Item(8, 17).Value = 14487.66779312379
$ws.Cells.Item(8, 18).Value = 130389.0101381141
$ws.Cells.Item(8, 19).Value = 0.1472592551478322
$ws.Cells.Item(8, 20).Value = 0.1472592551478323
$ws.Cells.Item(9, 9).Value = 0.7769829249672668
$ws.Cells.Item(9, 10).Value = 0.776982924967267
$ws.Cells.Item(9, 13).Value = 54.762539
$ws.Cells.Item(9, 14).Value = 164.287617
$ws.Cells.Item(9, 15).Value = 0.3866872460670236
$ws.Cells.Item(9, 16).Value = 0.3866872460670236
$ws.Cells.Item(9, 17).Value = 29558.82745929868
$ws.Cells.Item(9, 18).Value = 266029.4471336881
$ws.Cells.Item(9, 19).Value = 0.3004493874966933
$ws.Cells.Item(9, 20).Value = 0.3004493874966933
$ws.Cells.Item(10, 7).Value = 53.798087
$ws.Cells.Item(10, 8).Value = 161.394261
$ws.Cells.Item(10, 9).Value = 0.07744166482137986
$ws.Cells.Item(10, 10).Value = 0.07744166482137986
$ws.Cells.Item(10, 13).Value = 9.423852333333334
$ws.Cells.Item(10, 14).Value = 28.271557
$ws.Cells.Item(10, 15).Value = 0.06654336290212845
$ws.Cells.Item(10, 16).Value = 0.06654336290212845
$ws.Cells.Item(10, 17).Value = 506.9852277038198
$ws.Cells.Item(10, 18).Value = 4562.867049334377
$ws.Cells.Item(10, 19).Value = 0.005153228805954074
$ws.Cells.Item(10, 20).Value = 0.005153228805954074
$ws.Cells.Item(11, 7).Value = 53.798087
$ws.Cells.Item(11, 8).Value = 161.394261
$ws.Cells.Item(11, 9).Value = 0.07744166482137986
$ws.Cells.Item(11, 10).Value = 0.07744166482137986
$ws.Cells.Item(11, 15).Value = 0.3572423751649123
$ws.Cells.Item(11, 16).Value = 0.3572423751649123
$ws.Cells.Item(11, 17).Value = 2721.78319549046
$ws.Cells.Item(11, 18).Value = 24496.04875941413
$ws.Cells.Item(11, 19).Value = 0.02766544427751478
$ws.Cells.Item(11, 20).Value = 0.02766544427751478
$ws.Cells.Item(12, 7).Value = 53.798087
$ws.Cells.Item(12, 8).Value = 161.394261
$ws.Cells.Item(12, 9).Value = 0.07744166482137986
$ws.Cells.Item(12, 10).Value = 0.07744166482137986
$ws.Cells.Item(12, 13).Value = 26.84076266666667
$ws.Cells.Item(12, 14).Value = 80.522288
$ws.Cells.Item(12, 15).Value = 0.1895270158659356
$ws.Cells.Item(12, 16).Value = 0.1895270158659356
$ws.Cells.Item(12, 17).Value = 1443.981685087685
$ws.Cells.Item(12, 18).Value = 12995.83516578917
$ws.Cells.Item(12, 19).Value = 0.01467728763728613
$ws.Cells.Item(12, 20).Value = 0.01467728763728613
$ws.Cells.Item(13, 7).Value = 53.798087
$ws.Cells.Item(13, 8).Value = 161.394261
$ws.Cells.Item(13, 9).Value = 0.07744166482137986
$ws.Cells.Item(13, 10).Value = 0.07744166482137986
$ws.Cells.Item(13, 13).Value = 54.762539
$ws.Cells.Item(13, 14).Value = 164.287617
$ws.Cells.Item(13, 15).Value = 0.3866872460670236
$ws.Cells.Item(13, 16).Value = 0.3866872460670236
$ws.Cells.Item(13, 17).Value = 2946.119837462893
$ws.Cells.Item(13, 18).Value = 26515.07853716603
$ws.Cells.Item(13, 19).Value = 0.02994570410062488
$ws.Cells.Item(13, 20).Value = 0.02994570410062488
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 0.2555676666666666
$ws.Cells.Item(14, 8).Value = 0.7667029999999999
$ws.Cells.Item(14, 9).Value = 0.0003678864191059829
$ws.Cells.Item(14, 10).Value = 0.000367886419105983
$ws.Cells.Item(14, 13).Value = 9.423852333333334
$ws.Cells.Item(14, 14).Value = 28.271557
$ws.Cells.Item(14, 15).Value = 0.06654336290212845
$ws.Cells.Item(14, 16).Value = 0.06654336290212845
$ws.Cells.Item(14, 17).Value = 2.408431951841222
$ws.Cells.Item(14, 18).Value = 21.675887566571
$ws.Cells.Item(14, 19).Value = 0.00002448039949333394
$ws.Cells.Item(14, 20).Value = 0.00002448039949333395
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 0.2555676666666666
$ws.Cells.Item(15, 8).Value = 0.7667029999999999
$ws.Cells.Item(15, 9).Value = 0.0003678864191059829
$ws.Cells.Item(15, 10).Value = 0.000367886419105983
$ws.Cells.Item(15, 15).Value = 0.3572423751649123
$ws.Cells.Item(15, 16).Value = 0.3572423751649123
$ws.Cells.Item(15, 17).Value = 12.92982370254244
$ws.Cells.Item(15, 18).Value = 116.368413322882
$ws.Cells.Item(15, 19).Value = 0.0001314246181523357
$ws.Cells.Item(15, 20).Value = 0.0001314246181523357
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 0.2555676666666666
$ws.Cells.Item(16, 8).Value = 0.7667029999999999
$ws.Cells.Item(16, 9).Value = 0.0003678864191059829
$ws.Cells.Item(16, 10).Value = 0.000367886419105983
$ws.Cells.Item(16, 13).Value = 26.84076266666667
$ws.Cells.Item(16, 14).Value = 80.522288
$ws.Cells.Item(16, 15).Value = 0.1895270158659356
$ws.Cells.Item(16, 16).Value = 0.1895270158659356
$ws.Cells.Item(16, 17).Value = 6.859631086273777
$ws.Cells.Item(16, 18).Value = 61.73667977646399
$ws.Cells.Item(16, 19).Value = 0.00006972441519076184
$ws.Cells.Item(16, 20).Value = 0.00006972441519076187
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 0.2555676666666666
$ws.Cells.Item(17, 8).Value = 0.7667029999999999
$ws.Cells.Item(17, 9).Value = 0.0003678864191059829
$ws.Cells.Item(17, 10).Value = 0.000367886419105983
$ws.Cells.Item(17, 13).Value = 54.762539
$ws.Cells.Item(17, 14).Value = 164.287617
$ws.Cells.Item(17, 15).Value = 0.3866872460670236
$ws.Cells.Item(17, 16).Value = 0.3866872460670236
$ws.Cells.Item(17, 17).Value = 13.99553431297233
$ws.Cells.Item(17, 18).Value = 125.959808816751
$ws.Cells.Item(17, 19).Value = 0.0001422569862695514
$ws.Cells.Item(17, 20).Value = 0.0001422569862695514
